$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 1874
$ws.Cells.Item(6, 6).Value = 722
$ws.Cells.Item(7, 6).Value = 103
$ws.Cells.Item(8, 6).Value = 496
$ws.Cells.Item(9, 6).Value = 863
$ws.Cells.Item(10, 6).Value = 1577
$ws.Cells.Item(11, 6).Value = 1260
$ws.Cells.Item(12, 6).Value = 1505
$ws.Cells.Item(13, 6).Value = 52
$ws.Cells.Item(14, 6).Value = 1458
$ws.Cells.Item(15, 6).Value = 335
$ws.Cells.Item(16, 6).Value = 1677
$ws.Cells.Item(18, 6).Value = 1087
$ws.Cells.Item(19, 6).Value = 360
$ws.Cells.Item(20, 6).Value = 52
$ws.Cells.Item(22, 6).Value = 1673
$ws.Cells.Item(23, 6).Value = 206
$ws.Cells.Item(24, 6).Value = 814
$ws.Cells.Item(25, 6).Value = 556
$ws.Cells.Item(26, 6).Value = 1187
$ws.Cells.Item(27, 6).Value = 305678
$ws.Cells.Item(28, 6).Value = 1041
$ws.Cells.Item(29, 6).Value = 71
$ws.Cells.Item(32, 6).Value = 1133
$ws.Cells.Item(33, 6).Value = 901
$ws.Cells.Item(34, 6).Value = 2
$ws.Cells.Item(35, 6).Value = 1124
$ws.Cells.Item(36, 6).Value = 73
$ws.Cells.Item(37, 6).Value = 249
$ws.Cells.Item(38, 6).Value = 69
$ws.Cells.Item(39, 6).Value = 870
$ws.Cells.Item(40, 6).Value = 1668
$ws.Cells.Item(41, 6).Value = 6
$ws.Cells.Item(42, 6).Value = 114
$ws.Cells.Item(43, 6).Value = 8
$ws.Cells.Item(44, 6).Value = 84
$ws.Cells.Item(45, 6).Value = 823
$ws.Cells.Item(46, 6).Value = 91
$ws.Cells.Item(48, 6).Value = 115

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(5, 6).Value = 181
$ws.Cells.Item(9, 6).Value = 2576
$ws.Cells.Item(10, 6).Value = 1210
$ws.Cells.Item(11, 6).Value = 407
$ws.Cells.Item(13, 6).Value = 246
$ws.Cells.Item(14, 6).Value = 33
$ws.Cells.Item(15, 6).Value = 73
$ws.Cells.Item(18, 6).Value = 456
$ws.Cells.Item(21, 6).Value = 310
$ws.Cells.Item(22, 6).Value = 81957
$ws.Cells.Item(24, 6).Value = 2
$ws.Cells.Item(27, 6).Value = 189
$ws.Cells.Item(28, 6).Value = 247
$ws.Cells.Item(30, 6).Value = 210
$ws.Cells.Item(31, 6).Value = 59
$ws.Cells.Item(33, 6).Value = 56
$ws.Cells.Item(35, 6).Value = 178
$ws.Cells.Item(39, 6).Value = 59
$ws.Cells.Item(40, 6).Value = 59
$ws.Cells.Item(43, 6).Value = 63

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 118
$ws.Cells.Item(4, 6).Value = 261
$ws.Cells.Item(5, 6).Value = 2875
$ws.Cells.Item(6, 6).Value = 4630
$ws.Cells.Item(9, 6).Value = 570
$ws.Cells.Item(10, 6).Value = 724
$ws.Cells.Item(11, 6).Value = 460
$ws.Cells.Item(12, 6).Value = 334
$ws.Cells.Item(13, 6).Value = 1050
$ws.Cells.Item(14, 6).Value = 277
$ws.Cells.Item(15, 6).Value = 647

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 1874
$ws.Cells.Item(3, 6).Value = 261
$ws.Cells.Item(4, 6).Value = 4630
$ws.Cells.Item(5, 6).Value = 724
$ws.Cells.Item(7, 6).Value = 334
$ws.Cells.Item(8, 6).Value = 334
$ws.Cells.Item(9, 6).Value = 1051
$ws.Cells.Item(10, 6).Value = 1051
$ws.Cells.Item(11, 6).Value = 496
$ws.Cells.Item(12, 6).Value = 863
$ws.Cells.Item(13, 6).Value = 2576
$ws.Cells.Item(14, 6).Value = 1210
$ws.Cells.Item(15, 6).Value = 1577
$ws.Cells.Item(16, 6).Value = 1260
$ws.Cells.Item(17, 6).Value = 1505
$ws.Cells.Item(18, 6).Value = 52
$ws.Cells.Item(19, 6).Value = 1458
$ws.Cells.Item(20, 6).Value = 246
$ws.Cells.Item(21, 6).Value = 335
$ws.Cells.Item(22, 6).Value = 73
$ws.Cells.Item(23, 6).Value = 1677
$ws.Cells.Item(25, 6).Value = 1087
$ws.Cells.Item(26, 6).Value = 360
$ws.Cells.Item(27, 6).Value = 647
$ws.Cells.Item(28, 6).Value = 647
$ws.Cells.Item(29, 6).Value = 456
$ws.Cells.Item(30, 6).Value = 1673
$ws.Cells.Item(32, 6).Value = 206
$ws.Cells.Item(33, 6).Value = 814
$ws.Cells.Item(34, 6).Value = 556
$ws.Cells.Item(35, 6).Value = 1187
$ws.Cells.Item(36, 6).Value = 310
$ws.Cells.Item(37, 6).Value = 1041
$ws.Cells.Item(38, 6).Value = 71
$ws.Cells.Item(40, 6).Value = 1133
$ws.Cells.Item(41, 6).Value = 901
$ws.Cells.Item(42, 6).Value = 1124
$ws.Cells.Item(44, 6).Value = 249
$ws.Cells.Item(45, 6).Value = 870
$ws.Cells.Item(46, 6).Value = 59
$ws.Cells.Item(47, 6).Value = 1668
$ws.Cells.Item(48, 6).Value = 114
$ws.Cells.Item(49, 6).Value = 84
$ws.Cells.Item(50, 6).Value = 823
$ws.Cells.Item(51, 6).Value = 59
